# This script updates specific F (想去人数 / want-to-go count) and one G (最低票价 / min price)
# cell values in each worksheet of the workbook, per the commit's regenerated data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 476
$ws.Range("F6").Value = 837
$ws.Range("F7").Value = 257
$ws.Range("F8").Value = 1232
$ws.Range("F9").Value = 355
$ws.Range("F11").Value = 888
$ws.Range("F18").Value = 2964
$ws.Range("F19").Value = 2632
$ws.Range("F26").Value = 5330
$ws.Range("F27").Value = 593
$ws.Range("F28").Value = 992
$ws.Range("F31").Value = 333
$ws.Range("F32").Value = 1113
$ws.Range("F33").Value = 71
$ws.Range("F34").Value = 59
$ws.Range("F35").Value = 294

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1142
$ws.Range("F5").Value = 7
$ws.Range("F10").Value = 35
$ws.Range("F13").Value = 8
$ws.Range("F21").Value = 3
$ws.Range("F24").Value = 320
$ws.Range("F26").Value = 3963
$ws.Range("F29").Value = 22
$ws.Range("F33").Value = 170
$ws.Range("F36").Value = 8

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2486
$ws.Range("F6").Value = 1058
$ws.Range("F9").Value = 1346

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2486
$ws.Range("F7").Value = 1058
$ws.Range("F8").Value = 1346
$ws.Range("F12").Value = 476
$ws.Range("F13").Value = 837
$ws.Range("F14").Value = 257
$ws.Range("F15").Value = 1232
$ws.Range("F16").Value = 355
$ws.Range("F17").Value = 888
$ws.Range("F19").Value = 1142
$ws.Range("F20").Value = 1142
$ws.Range("F25").Value = 2964
$ws.Range("F26").Value = 2632
$ws.Range("F29").Value = 35
$ws.Range("F32").Value = 5330
$ws.Range("F33").Value = 593
$ws.Range("F34").Value = 992
$ws.Range("F39").Value = 333
$ws.Range("F42").Value = 3
$ws.Range("F44").Value = 320
$ws.Range("F45").Value = 320
$ws.Range("F46").Value = 1113
$ws.Range("G48").Value = 380
$ws.Range("F49").Value = 170
$ws.Range("F50").Value = 59
$ws.Range("F51").Value = 294
